$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2066.5
$ws.Range("I2").Value = 1133
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 1133
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -1020
$ws.Range("N2").Value = -3226

$ws.Range("H43").Value = 3343.0625
$ws.Range("I43").Value = 2331.6667
$ws.Range("J43").Value = 3576.4614
$ws.Range("K43").Value = 2331.6667
$ws.Range("L43").Value = 3576.4614
$ws.Range("M43").Value = -2262.6667
$ws.Range("N43").Value = -3714.4614

$ws.Range("H126").Value = 74999.5
$ws.Range("J126").Value = 74999.5
$ws.Range("L126").Value = 74999.5
$ws.Range("N126").Value = -84879.5

$ws.Range("H137").Value = 28171.129
$ws.Range("I137").Value = 33493.625
$ws.Range("J137").Value = 3839.7144
$ws.Range("K137").Value = 100480.875
$ws.Range("L137").Value = 11519.1432
$ws.Range("M137").Value = -97930.875
$ws.Range("N137").Value = -16619.1432

$ws.Range("H138").Value = 2487.2183
$ws.Range("J138").Value = 2787.5186
$ws.Range("L138").Value = 8362.5558
$ws.Range("N138").Value = -18642.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 25453.023
$ws.Range("I132").Value = 33960.812
$ws.Range("K132").Value = 101882.436
$ws.Range("M132").Value = -99352.43599999999

$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1842.7441
$ws.Range("I94").Value = 1613.9678
$ws.Range("K94").Value = 1613.9678
$ws.Range("M94").Value = -1162.9678

$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1816.7241
$ws.Range("I31").Value = 1085.2858
$ws.Range("J31").Value = 3736.75
$ws.Range("K31").Value = 1085.2858
$ws.Range("L31").Value = 3736.75
$ws.Range("M31").Value = -790.2858000000001
$ws.Range("N31").Value = -4326.75

$ws.Range("H34").Value = 1816.7241
$ws.Range("I34").Value = 1085.2858
$ws.Range("J34").Value = 3736.75
$ws.Range("K34").Value = 1085.2858
$ws.Range("L34").Value = 3736.75
$ws.Range("M34").Value = -883.2858000000001
$ws.Range("N34").Value = -4140.75

$ws.Range("H86").Value = 4926.7144
$ws.Range("J86").Value = 4944.3335
$ws.Range("L86").Value = 4944.3335
$ws.Range("N86").Value = -7190.3335

$ws.Range("H89").Value = 4926.7144
$ws.Range("J89").Value = 4944.3335
$ws.Range("L89").Value = 24721.6675
$ws.Range("N89").Value = -35953.6675

$ws.Range("H132").Value = 2602.7778
$ws.Range("I132").Value = 2527.4666
$ws.Range("J132").Value = 2979.3333
$ws.Range("K132").Value = 7582.399800000001
$ws.Range("L132").Value = 8937.999899999999
$ws.Range("M132").Value = -5052.399800000001
$ws.Range("N132").Value = -13997.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8004
$ws.Range("I3").Value = 5005
$ws.Range("K3").Value = 15015
$ws.Range("M3").Value = -14903

$ws.Range("H56").Value = 7385.5
$ws.Range("I56").Value = 7385.5
$ws.Range("K56").Value = 7385.5
$ws.Range("M56").Value = -6855.5

$ws.Range("H81").Value = 3256
$ws.Range("I81").Value = 3256
$ws.Range("K81").Value = 9768
$ws.Range("M81").Value = -8645

$ws.Range("H84").Value = 3256
$ws.Range("I84").Value = 3256
$ws.Range("K84").Value = 29304
$ws.Range("M84").Value = -23688

$ws.Range("H108").Value = 1738
$ws.Range("I108").Value = 1738
$ws.Range("K108").Value = 5214
$ws.Range("M108").Value = -2334

$ws.Range("H114").Value = 819.1429000000001
$ws.Range("I114").Value = 701.625
$ws.Range("J114").Value = 975.8333
$ws.Range("K114").Value = 2104.875
$ws.Range("L114").Value = 2927.4999
$ws.Range("M114").Value = 1149.125
$ws.Range("N114").Value = -9435.499899999999

$ws.Range("H121").Value = 287.4
$ws.Range("I121").Value = 324
$ws.Range("J121").Value = 202
$ws.Range("K121").Value = 972
$ws.Range("L121").Value = 606
$ws.Range("M121").Value = 338
$ws.Range("N121").Value = -3226

$ws.Range("H136").Value = 4102
$ws.Range("I136").Value = 4102
$ws.Range("K136").Value = 12306
$ws.Range("M136").Value = -7206

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3714.5
$ws.Range("I80").Value = 3462.6667
$ws.Range("J80").Value = 3966.3333
$ws.Range("K80").Value = 3462.6667
$ws.Range("L80").Value = 3966.3333
$ws.Range("M80").Value = -2464.6667
$ws.Range("N80").Value = -5962.3333

$ws.Range("H83").Value = 3714.5
$ws.Range("I83").Value = 3462.6667
$ws.Range("J83").Value = 3966.3333
$ws.Range("K83").Value = 17313.3335
$ws.Range("L83").Value = 19831.6665
$ws.Range("M83").Value = -12321.3335
$ws.Range("N83").Value = -29815.6665

$ws.Range("H102").Value = 3240.182
$ws.Range("I102").Value = 3225.9333
$ws.Range("K102").Value = 3225.9333
$ws.Range("M102").Value = -1603.9333

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H126").Value = 6477.4546
$ws.Range("J126").Value = 7428.1113
$ws.Range("L126").Value = 22284.3339
$ws.Range("N126").Value = -27224.3339

$ws.Range("H140").Value = 89999
$ws.Range("J140").Value = 89999
$ws.Range("L140").Value = 89999
$ws.Range("N140").Value = -100359

$ws.Range("H141").Value = 69999
$ws.Range("J141").Value = 69999
$ws.Range("L141").Value = 69999
$ws.Range("N141").Value = -80359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 14476.5
$ws.Range("I46").Value = 19729.924
$ws.Range("J46").Value = 6888.222
$ws.Range("K46").Value = 19729.924
$ws.Range("L46").Value = 6888.222
$ws.Range("M46").Value = -19541.924
$ws.Range("N46").Value = -7264.222

$ws.Range("H55").Value = 561.1667
$ws.Range("J55").Value = 908.3333
$ws.Range("L55").Value = 908.3333
$ws.Range("N55").Value = -1254.3333

$ws.Range("H61").Value = 2152.6667
$ws.Range("I61").Value = 2020.3572
$ws.Range("K61").Value = 2020.3572
$ws.Range("M61").Value = -1818.3572

$ws.Range("H113").Value = 2152.6667
$ws.Range("I113").Value = 2020.3572
$ws.Range("K113").Value = 2020.3572
$ws.Range("M113").Value = 149.6428000000001

$ws.Range("H136").Value = 1803
$ws.Range("J136").Value = 2398.5
$ws.Range("L136").Value = 7195.5
$ws.Range("N136").Value = -12295.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1786.6
$ws.Range("I113").Value = 1373.6
$ws.Range("J113").Value = 2199.6
$ws.Range("K113").Value = 4120.799999999999
$ws.Range("L113").Value = 6598.799999999999
$ws.Range("M113").Value = -1950.799999999999
$ws.Range("N113").Value = -10938.8

$ws.Range("H132").Value = 21947.674
$ws.Range("I132").Value = 22862.836
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 68588.508
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -66058.508
$ws.Range("N132").Value = -26060

$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139

$ws.Range("H136").Value = 3646.5
$ws.Range("I136").Value = 3646.5
$ws.Range("K136").Value = 10939.5
$ws.Range("M136").Value = -8389.5
